$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Förändrad" (C) column timestamp advanced by 2 days (2024-06-09 -> 2024-06-11)
# for every existing data row (rows 2-28).
$ws.Range("C2:C28").Value = 45454

# The last data row (29, case "A 23041-2024") was removed from the log.
$ws.Rows("29:29").Delete()

# Deleting the trailing row leaves row 28 (now the last row) without an
# explicit custom row height, so re-fit it to drop the stale ht/customHeight.
$ws.Rows("28:28").AutoFit()
